$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "28.466.86"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +0.00%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.824.69"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -0.12%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.32%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "316.25"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +0.30%  "
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +0.25%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5186"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3866"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -0.99%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.08278"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +8.11%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.122"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +1.25%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "41.85"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -0.09%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "6.371"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +1.29%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "21.13"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +0.20%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.004"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.485"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -1.38%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.829.01"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +0.31%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "93.92"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +0.82%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001120"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +3.39%  "
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -0.63%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.82"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +0.79%  "
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +0.28%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.069"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -1.20%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "28.512.92"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +0.06%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.50"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +3.28%  "
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -0.45%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "21.10"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +2.37%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "159.56"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +1.77%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.037.99"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.410"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +0.60%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "125.95"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +0.55%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.1100"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +1.65%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.096"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -2.84%  "
$ws.Range("B33").Value = "Hedera"
$ws.Range("C33").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.07610"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +8.08%  "
$ws.Range("B34").Value = "Filecoin"
$ws.Range("C34").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.733"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +0.90%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.683"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +0.59%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.2228"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -0.05%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02369"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +2.12%  "
$ws.Range("B38").Value = "Aptos"
$ws.Range("C38").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "12.10"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +7.80%  "
$ws.Range("B39").Value = "InternetComputer(DFINITY)"
$ws.Range("C39").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.254"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +2.30%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "8.762"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -2.13%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.6418"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +2.56%  "
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +1.02%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.401"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +0.24%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "13.72"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +2.39%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.6240"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +5.74%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.798"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +2.18%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "127.88"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +2.74%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.206"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +1.05%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.06973"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.080"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +1.25%  "
